$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "d"
$ws.Range("C2").Value = "d"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 0
$ws.Range("H2").Value = 54.7
